$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows of daily user-impact data to append below the existing table
# (rows 33-38): Date, # Total Users, # Logged-in Users, 0 Errors, 1 Errors,
# 2 Errors, 3-5 Errors, 6-10 Errors, >10 Errors
$newRows = @(
    @(45982, 5627, 4092, 3816, 197, 46, 27, 6, 0),
    @(45983, 5627, 56,   56,   0,   0,  0,  0, 0),
    @(45984, 5627, 42,   42,   0,   0,  0,  0, 0),
    @(45985, 5626, 4225, 3909, 227, 45, 39, 4, 1),
    @(45986, 5623, 4194, 3868, 245, 52, 24, 5, 0),
    @(45987, 5620, 3695, 3413, 215, 33, 30, 3, 1)
)

$startRow = 33
for ($i = 0; $i -lt $newRows.Length; $i++) {
    $r = $startRow + $i
    $rowValues = $newRows[$i]
    for ($c = 1; $c -le 9; $c++) {
        $ws.Cells.Item($r, $c).Value = $rowValues[$c - 1]
    }
}

$lastRow = $startRow + $newRows.Length - 1

# Copy formatting (number formats/styles) from the last previously-existing
# data row (32) down onto the newly added rows, matching how Excel extends
# formatting when a user fills new rows below existing data.
$ws.Range("A32:I32").Copy()
$ws.Range("A33:I$lastRow").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Select the newly added last row, matching Excel's behaviour of landing on
# the last entered row after data entry.
$ws.Range("A$lastRow`:I$lastRow").Select()
